# Auto-generated Excel COM-interop script to apply diff changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("U2").Value = 1.47
$ws.Range("W2").Value = 11
$ws.Range("G3").Value = 2.15
$ws.Range("H3").Value = 3.3
$ws.Range("I3").Value = 3.3
$ws.Range("J3").Value = 2.75
$ws.Range("L3").Value = 3.75
$ws.Range("M3").Value = 1.03
$ws.Range("O3").Value = 1.22
$ws.Range("Q3").Value = 1.95
$ws.Range("R3").Value = 1.9
$ws.Range("U3").Value = 1.69
$ws.Range("Y3").Value = 9
$ws.Range("AA3").Value = 17
$ws.Range("AH3").Value = 11
$ws.Range("AL3").Value = 26
$ws.Range("AX3").Value = 5.5
$ws.Range("BB3").Value = 81
$ws.Range("G4").Value = 5.25
$ws.Range("H4").Value = 4.2
$ws.Range("I4").Value = 1.53
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 2.5
$ws.Range("L4").Value = 2.05
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 17
$ws.Range("O4").Value = 1.17
$ws.Range("P4").Value = 5
$ws.Range("Q4").Value = 1.57
$ws.Range("R4").Value = 2.35
$ws.Range("S4").Value = 1.29
$ws.Range("T4").Value = 3.5
$ws.Range("U4").Value = 1.63
$ws.Range("V4").Value = 2.1
$ws.Range("W4").Value = 19
$ws.Range("X4").Value = 29
$ws.Range("Y4").Value = 17
$ws.Range("Z4").Value = 51
$ws.Range("AA4").Value = 41
$ws.Range("AB4").Value = 41
$ws.Range("AD4").Value = 8.5
$ws.Range("AE4").Value = 15
$ws.Range("AI4").Value = 8.5
$ws.Range("AK4").Value = 12
$ws.Range("AL4").Value = 12
$ws.Range("AN4").Value = 7
$ws.Range("AO4").Value = 26
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 81
$ws.Range("AT4").Value = 3.5
$ws.Range("AU4").Value = 8
$ws.Range("AX4").Value = 3.75
$ws.Range("AY4").Value = 7.5
$ws.Range("AZ4").Value = 15
$ws.Range("BA4").Value = 21
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("Q6").Value = 1.98
$ws.Range("R6").Value = 1.88
